$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# B5: "write code methods for same;" -> "write code methods listed below"
# ------------------------------------------------------------------
$ws.Range("B5").Value2 = "write code methods listed below"

# ------------------------------------------------------------------
# "methods in code.js to finish" TODO list (rows 15-22):
#   remove "self.isReferencedInWorkspace" entirely, re-order so the
#   three completed items (replaceMethod / replaceType / replaceProperty)
#   sit at the bottom of the list and get struck through, and bold+
#   underline the list header.
# ------------------------------------------------------------------
$ws.Range("A16").Value2 = "self.removeEvent"
$ws.Range("A17").Value2 = "self.renameEvent"
$ws.Range("A18").Value2 = "m_functionRemove_Type_Event"
$ws.Range("E18").Value2 = "see m_functionAdd_Type_Event"
$ws.Range("A19").Value2 = "self.replaceMethod"
$ws.Range("A20").Value2 = "self.replaceType"
$ws.Range("A21").Value2 = "self.replaceProperty"

# old rows 22 and 27 no longer hold these values after the re-layout
$ws.Range("A22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("A27").ClearContents()

# "Check out (...)" note moves up from row 27 to row 25
$ws.Range("A25").Value2 = "Check out (in code.js) self.renameTypeInActiveComic and self.renameMethodInActiveType to be sure everything's being done."

# New trailing note
$ws.Range("A28").Value2 = "Save project to DB"

# ------------------------------------------------------------------
# Formatting: bold+underline the header, strikethrough the done items
# ------------------------------------------------------------------
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Font.Underline = $true

$ws.Range("A19:A21").Font.Strikethrough = $true

# ------------------------------------------------------------------
# Selection, matching the saved workbook view
# ------------------------------------------------------------------
$ws.Range("A22").Select() | Out-Null
